# Edit script: reorders data rows 2-132 according to the target order of the
# "Beteckning" (column A) key, and updates column C ("Förändrad") from 46070 to 46072
# for every data row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$firstRow = 2
$lastRow = 132
$rowCount = $lastRow - $firstRow + 1

$targetOrder = @(
  "A 14442-2025",
  "A 2443-2023",
  "A 2440-2023",
  "A 28233-2022",
  "A 55111-2025",
  "A 2447-2023",
  "A 22268-2022",
  "A 3543-2024",
  "A 64855-2023",
  "A 15136-2022",
  "A 2438-2023",
  "A 33922-2022",
  "A 72415-2021",
  "A 58025-2021",
  "A 12978-2021",
  "A 30197-2022",
  "A 29106-2022",
  "A 22300-2022",
  "A 20681-2022",
  "A 29176-2021",
  "A 13411-2022",
  "A 468-2022",
  "A 22240-2022",
  "A 15138-2022",
  "A 28249-2021",
  "A 52777-2021",
  "A 29101-2022",
  "A 69614-2021",
  "A 69628-2021",
  "A 26345-2022",
  "A 21133-2022",
  "A 50622-2023",
  "A 39113-2024",
  "A 55429-2025",
  "A 38743-2025",
  "A 38739-2025",
  "A 59758-2024",
  "A 14409-2025",
  "A 56672-2025",
  "A 13156-2024",
  "A 40848-2025",
  "A 26073-2025",
  "A 50420-2023",
  "A 35199-2023",
  "A 11124-2025",
  "A 43043-2025",
  "A 27367-2025",
  "A 28053-2025",
  "A 28029-2025",
  "A 58073-2025",
  "A 20210-2023",
  "A 58086-2025",
  "A 45132-2025",
  "A 44778-2025",
  "A 44780-2025",
  "A 29394-2025",
  "A 45396-2025",
  "A 45363-2025",
  "A 30518-2021",
  "A 58042-2025",
  "A 45483-2025",
  "A 13100-2025",
  "A 28089-2022",
  "A 46596-2025",
  "A 47134-2025",
  "A 46957-2025",
  "A 47599-2025",
  "A 66104-2021",
  "A 57186-2022",
  "A 31687-2025",
  "A 57083-2022",
  "A 57120-2022",
  "A 57154-2022",
  "A 32093-2025",
  "A 36756-2022",
  "A 2450-2023",
  "A 66107-2021",
  "A 66008-2021",
  "A 32095-2025",
  "A 57178-2022",
  "A 57195-2022",
  "A 51900-2021",
  "A 23592-2022",
  "A 7272-2023",
  "A 37441-2022",
  "A 22223-2022",
  "A 50397-2023",
  "A 28103-2022",
  "A 7989-2024",
  "A 10410-2022",
  "A 50405-2023",
  "A 50410-2023",
  "A 47866-2024",
  "A 41836-2022",
  "A 2420-2026",
  "A 2335-2026",
  "A 61762-2025",
  "A 61766-2025",
  "A 61628-2025",
  "A 47398-2024",
  "A 62736-2025",
  "A 7429-2024",
  "A 62755-2025",
  "A 62769-2025",
  "A 62774-2025",
  "A 62747-2025",
  "A 62757-2025",
  "A 17071-2024",
  "A 4883-2024",
  "A 20570-2023",
  "A 7558-2026",
  "A 26568-2023",
  "A 20680-2022",
  "A 40709-2023",
  "A 24163-2024",
  "A 14067-2024",
  "A 5359-2023",
  "A 36415-2023",
  "A 42526-2024",
  "A 47944-2024",
  "A 28118-2024",
  "A 28139-2024",
  "A 7990-2024",
  "A 64866-2023",
  "A 16648-2023",
  "A 48941-2024",
  "A 790-2024",
  "A 28157-2022",
  "A 17861-2021",
  "A 11232-2025",
  "A 21836-2022"
)

$totalRows = $targetOrder.Count
Write-Host "Target order contains $totalRows entries; expecting $rowCount"

# --- Step 1: capture the current row number for each Beteckning value ---
$currentRowForKey = @{}
for ($r = $firstRow; $r -le $lastRow; $r++) {
    $key = $ws.Range("A$r").Value2
    $currentRowForKey[$key] = $r
}

# --- Step 2: snapshot ALL current row data (values for A:R, formulas for S:Z) ---
# before any writes happen, so source rows are not clobbered mid-move.
$snapshotValues = @{}   # row -> object[] (A..R, i.e. 18 columns)
$snapshotFormulas = @{} # row -> object[] (S..Z, i.e. 8 columns)

for ($r = $firstRow; $r -le $lastRow; $r++) {
    $snapshotValues[$r] = $ws.Range("A$r`:R$r").Value2
    $snapshotFormulas[$r] = $ws.Range("S$r`:Z$r").Formula
}

# --- Step 3: write each row's data into its new target position ---
for ($i = 0; $i -lt $totalRows; $i++) {
    $destRow = $firstRow + $i
    $key = $targetOrder[$i]
    $srcRow = $currentRowForKey[$key]

    $valRange = $ws.Range("A$destRow`:R$destRow")
    $valRange.Value2 = $snapshotValues[$srcRow]

    $formRange = $ws.Range("S$destRow`:Z$destRow")
    $formRange.Formula = $snapshotFormulas[$srcRow]
}

# --- Step 4: update column C ("Förändrad") to 46072 for every data row ---
for ($r = $firstRow; $r -le $lastRow; $r++) {
    $ws.Range("C$r").Value2 = 46072
}

Write-Host "Row reorder and column C update complete."
